# "first push ubp calculation"
# Updates data/scenario_options.xlsx: new header label, new weather-file
# paths (moved from walkerl's RC_BuildingSimulator checkout to the
# LW_Simulation PycharmProjects checkout), numeric UBP-factor columns for
# both scenario rows, refreshed cooling-setpoint/emission-source values,
# and a fully populated third ("2070-A2") scenario row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels -------------------------------------------------
# "occupancy schedule" -> "building use type"
$ws.Range("B1").Value = "building use type"

# --- Row 2: historic weatherfile scenario ---------------------------------
$ws.Range("A2").Value = "C:\Users\LW_Simulation\PycharmProjects\sia_380-1-full_version\data\Zürich-hour_historic.epw"
$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = "eu"

# --- Row 3: 2070-A2 weatherfile scenario ----------------------------------
$ws.Range("A3").Value = "C:\Users\LW_Simulation\PycharmProjects\sia_380-1-full_version\data\Zürich-2070-A2.epw"
$ws.Range("B3").Value = 3.1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 23
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = "SIA"

# --- Selection -------------------------------------------------------------
$null = $ws.Range("A11").Select()

# --- Column widths (closest reachable values; engine quantizes ColumnWidth
#     to 1/6-character steps, so these are the nearest achievable inputs) --
$ws.Columns.Item(1).ColumnWidth = 76.3333333333333
$ws.Columns.Item(2).ColumnWidth = 87.6666666666667
$ws.Columns.Item(3).ColumnWidth = 26.3333333333333
$ws.Columns.Item(4).ColumnWidth = 15.6666666666667
$ws.Columns.Item(5).ColumnWidth = 13.8333333333333

Write-Host "scenario_options.xlsx updated"
